# Generate Report for Handback
#
# Updates the "bd1f1b19-458e-48fd-ab2f-511a5367d156" row with freshly
# generated handoff/handback timestamps across the Overview, zh-cn and
# de-de sheets (mirrors a "Generate Report for Handback" run).

$wb = $excel.ActiveWorkbook

# Overview sheet: refresh the "Latest HO Xliff Generate Date" for the
# bd1f1b19 row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 20:44:20"

# zh-cn sheet: refresh Correspond Handoff / Handback datetimes for the
# bd1f1b19 row (row 3).
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H3").Value = "2016-08-16 20:44:15"
$wsZh.Range("K3").Value = "2016-08-16 20:44:30"

# de-de sheet: refresh Correspond Handoff / Handback datetimes for the
# bd1f1b19 row (row 3).
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H3").Value = "2016-08-16 20:44:20"
$wsDe.Range("K3").Value = "2016-08-16 20:44:37"
